# station_load_graph.xlsx update
# - Rows 2-17 (Hour 1-16): refresh the load readings for each station.
# - Rows 18-23 (Hour 17-22) and their data are removed (cleared back to empty rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New readings, keyed by row number: B,C,D,E,F,G,H,I (Lapaz, Jaro, Mandurriao, Molo,
# Diversion, Mobile SS 1, Mobile SS 2, Megaworld)
$rows = @{
  2  = @(18849, 4993,  8809, 0, 5861,  0, 14991, 18553)
  3  = @(18175, 4852,  8549, 0, 5503,  0, 14425, 17804)
  4  = @(17755, 4595,  8260, 0, 5166,  0, 12730, 17106)
  5  = @(16961, 4508,  8095, 0, 5128,  0, 13244, 16333)
  6  = @(17567, 4568,  8089, 0, 5068,  0, 12963, 16881)
  7  = @(18504, 5079,  8845, 0, 6427,  0, 14292, 17918)
  8  = @(18075, 4252,  9078, 0, 7144,  0, 14101, 17399)
  9  = @(20157, 4483, 10453, 0, 9822,  0, 16102, 19467)
  10 = @(25877, 4973, 12810, 0, 12116, 0, 17894, 21055)
  11 = @(27231, 5261, 16678, 0, 13582, 0, 18314, 23057)
  12 = @(28382, 5591, 17993, 0, 14269, 0, 19673, 24404)
  13 = @(29040, 5776, 18247, 0, 14444, 0, 19647, 24167)
  14 = @(27988, 5392, 18188, 0, 14483, 0, 18997, 24211)
  15 = @(29908, 5917, 18554, 0, 14599, 0, 20353, 24098)
  16 = @(29842, 10569, 18693, 0, 14775, 0, 20463, 24821)
  17 = @(29762, 1242, 18487, 0, 0,     0, 20385, 0)
}

$columns = @("B", "C", "D", "E", "F", "G", "H", "I")

foreach ($r in 2..17) {
  $values = $rows[$r]
  for ($i = 0; $i -lt $columns.Length; $i++) {
    $ws.Range("$($columns[$i])$r").Value = $values[$i]
  }
}

# Hours 17-22 (rows 18-23) were removed from the source data entirely.
$ws.Range("A18:I23").ClearContents()
